{"js": "// Update the two-digit division worksheet: replace each \"A\u00f7B=\" problem\n// text with its new value, per the target diff.\n//\n// Strategy: some new values collide with old values elsewhere in the\n// document (e.g. \"19\u00f74=\" becomes \"56\u00f75=\", while the original \"56\u00f75=\"\n// becomes \"21\u00f75=\"). To avoid a later replacement accidentally re-matching\n// text that an earlier replacement just inserted, we first run ALL the\n// searches against the pristine document and capture their Range objects,\n// then apply every insertText(...) afterwards.\n\nconst mapping = [\n  [\"16\u00f73=\", \"46\u00f79=\"],\n  [\"56\u00f77=\", \"37\u00f78=\"],\n  [\"54\u00f72=\", \"44\u00f77=\"],\n  [\"55\u00f74=\", \"48\u00f72=\"],\n  [\"68\u00f76=\", \"22\u00f73=\"],\n  [\"47\u00f75=\", \"39\u00f79=\"],\n  [\"25\u00f74=\", \"55\u00f72=\"],\n  [\"22\u00f77=\", \"71\u00f78=\"],\n  [\"26\u00f79=\", \"19\u00f79=\"],\n  [\"68\u00f74=\", \"37\u00f73=\"],\n  [\"19\u00f74=\", \"56\u00f75=\"],\n  [\"12\u00f78=\", \"83\u00f74=\"],\n  [\"93\u00f75=\", \"68\u00f73=\"],\n  [\"69\u00f76=\", \"65\u00f75=\"],\n  [\"78\u00f75=\", \"28\u00f78=\"],\n  [\"83\u00f76=\", \"70\u00f72=\"],\n  [\"77\u00f76=\", \"64\u00f75=\"],\n  [\"56\u00f75=\", \"21\u00f75=\"],\n  [\"58\u00f79=\", \"99\u00f78=\"],\n  [\"86\u00f74=\", \"29\u00f76=\"],\n  [\"50\u00f79=\", \"27\u00f78=\"],\n  [\"48\u00f77=\", \"40\u00f76=\"],\n  [\"22\u00f79=\", \"16\u00f75=\"],\n  [\"76\u00f77=\", \"54\u00f75=\"],\n  [\"70\u00f76=\", \"89\u00f73=\"],\n];\n\n// Phase 1: search for every old value in the untouched document and keep\n// the resulting SearchResult collections around.\nconst searchResults = mapping.map(([oldText]) =>\n  context.document.body.search(oldText, { matchCase: true, matchWholeWord: false })\n);\nsearchResults.forEach((res) => res.load(\"items\"));\nawait context.sync();\n\nsearchResults.forEach((res, i) => {\n  const [oldText] = mapping[i];\n  if (res.items.length !== 1) {\n    throw new Error(`Expected exactly one match for \"${oldText}\", found ${res.items.length}`);\n  }\n});\n\n// Phase 2: replace the text in each captured range with its new value.\nsearchResults.forEach((res, i) => {\n  const [, newText] = mapping[i];\n  res.items[0].insertText(newText, Word.InsertLocation.replace);\n});\nawait context.sync();\n", "ps1": "# Update the two-digit division worksheet: replace each \"A\u00f7B=\" problem\n# text with its new value, per the target diff.\n#\n# The worksheet is a single 5-column table; only every 4th row (1, 5, 9,\n# 13, 17) actually holds a problem, the rows in between are blank spacer\n# rows. We address each cell directly by (row, column) instead of using\n# Find/Replace so that a replacement value which happens to equal another\n# cell's original text (e.g. \"19\u00f74=\" becomes \"56\u00f75=\", while the original\n# \"56\u00f75=\" cell becomes \"21\u00f75=\") can never be re-matched by a later step.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$contentRows = @(1, 5, 9, 13, 17)\n\n# New text for each cell, in row-major order matching $contentRows x 5 columns.\n$newValues = @(\n    \"46\u00f79=\", \"37\u00f78=\", \"44\u00f77=\", \"48\u00f72=\", \"22\u00f73=\",\n    \"39\u00f79=\", \"55\u00f72=\", \"71\u00f78=\", \"19\u00f79=\", \"37\u00f73=\",\n    \"56\u00f75=\", \"83\u00f74=\", \"68\u00f73=\", \"65\u00f75=\", \"28\u00f78=\",\n    \"70\u00f72=\", \"64\u00f75=\", \"21\u00f75=\", \"99\u00f78=\", \"29\u00f76=\",\n    \"27\u00f78=\", \"40\u00f76=\", \"16\u00f75=\", \"54\u00f75=\", \"89\u00f73=\"\n)\n\n$i = 0\nforeach ($row in $contentRows) {\n    for ($col = 1; $col -le 5; $col++) {\n        $cell = $t.Cell($row, $col)\n        $r = $cell.Range\n        # Drop the trailing end-of-cell marker so we only overwrite the\n        # visible text, not the cell/paragraph mark.\n        $r.End = $r.End - 1\n        $r.Text = $newValues[$i]\n        $i = $i + 1\n    }\n}\n\nWrite-Output \"done\"\n"}
